# ---------------------------------------------------------------------------
# Edit script: rewrite Dheeraj Chand resume (short / modern_tech ATS variant)
# to match the "comprehensive inheritance system" commit.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

# Replace the full text of the (unique) paragraph that contains $findText
# with $newText, in place (keeps the paragraph's own formatting/style).
function Replace-ParaText($findText, $newText) {
    $r = $d.Content
    $found = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        Write-Output "REPLACE FAILED (not found): $findText"
    }
}

# Delete the paragraph range starting at the paragraph containing $startText
# through (and including) the paragraph containing $endText.
function Delete-ParaRange($startText, $endText) {
    $r1 = $d.Content
    $r1.Find.Execute($startText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $r1.Find.Found) {
        Write-Output "DELETE-RANGE START NOT FOUND: $startText"
        return
    }
    $startPos = $r1.Paragraphs(1).Range.Start

    $r2 = $d.Content
    $r2.Find.Execute($endText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $r2.Find.Found) {
        Write-Output "DELETE-RANGE END NOT FOUND: $endText"
        return
    }
    $endPos = $r2.Paragraphs(1).Range.End

    $delRange = $d.Range($startPos, $endPos)
    $delRange.Delete()
}

# Insert a new block of paragraphs right after the (unique) paragraph that
# contains $anchorText. $lines is an ordered array of strings, one per new
# paragraph. If $headingStyle is non-empty, it is applied to the first
# paragraph of the block (the job/project title line).
function Insert-Block($anchorText, $lines, $headingStyle) {
    $r = $d.Content
    $r.Find.Execute($anchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $r.Find.Found) {
        Write-Output "INSERT ANCHOR NOT FOUND: $anchorText"
        return
    }
    $anchorPara = $r.Paragraphs(1)
    $anchorPara.Range.InsertParagraphAfter()
    $newPara = $anchorPara.Next()
    $newPara.Range.Text = [string]::Join("`r", $lines)
    if ($headingStyle) {
        $newPara.Style = $headingStyle
    }
}

# ---------------------------------------------------------------------------
# 1. Remove the centered contact-info paragraph under the name.
# ---------------------------------------------------------------------------

$r = $d.Content
$r.Find.Execute('+1 (512) 555-0123 | dheeraj@dheerajchand.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $r.Paragraphs(1).Range.Delete()
} else {
    Write-Output "CONTACT LINE NOT FOUND"
}

# ---------------------------------------------------------------------------
# 2. Professional summary rewrite.
# ---------------------------------------------------------------------------

Replace-ParaText 'Experienced data scientist and software engineer with 15+ years of expertise in geospatial analysis, demographic research, and political data. Proven track record of building scalable systems, conducting complex analyses, and delivering actionable insights for campaigns, organizations, and government agencies.' 'Senior data scientist and software engineer specializing in geospatial machine learning and large-scale demographic analysis. Developed algorithms that improved demographic classification accuracy from 23% to 64%, processed data across 178,000+ precincts, and built platforms serving thousands of analysts nationwide.'

# ---------------------------------------------------------------------------
# 3. Core competencies: drop the tag line, leave an empty paragraph.
# ---------------------------------------------------------------------------

$r = $d.Content
$r.Find.Execute('CODE • COMPUTE • INTERACT • MEASURE • PLATFORMS • TRACK', $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($r.Find.Found) {
    $r.Paragraphs(1).Range.Text = ""
} else {
    Write-Output "CORE COMPETENCIES LINE NOT FOUND"
}

# ---------------------------------------------------------------------------
# 4. Professional experience - Siege Analytics (kept position, new copy).
# ---------------------------------------------------------------------------

Replace-ParaText 'Partner - Siege Analytics (Austin, TX) | 2020 - Present' 'Partner - Siege Analytics (Austin, TX) | 2005 - Present'
Replace-ParaText 'Data Science & Political Analytics' 'Data, Technology and Strategy Consulting'
Replace-ParaText '• Uncovered decades of demographic miscoding in voter files, discovering 2.7M previously mischaracterized Democratic voters' '• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from 23% to 64%'
Replace-ParaText '• Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States' '• Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration'
Replace-ParaText '• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct redistricting analysis' '• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%'

# ---------------------------------------------------------------------------
# 5. Professional experience - remove the old 6 jobs after Siege Analytics
#    (Lake Research Partners Sr Data Scientist .. Feldman Group), replace
#    with the 7 new roles, keeping the document's "Siege Analytics is
#    strongest / first" ordering intact.
# ---------------------------------------------------------------------------

Delete-ParaRange 'Senior Data Scientist - Lake Research Partners (Washington, DC) | 2018 - 2020' '• Trained staff on PHP/MySQL for data analysis and reporting systems'

Insert-Block '• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%, reducing polling error margins from ±4.2% to ±2.1%' @(
    'Data Products Manager - Helm/Murmuration (Austin, TX) | June 2021 - May 2023',
    'Civic Graph & Civic Pulse Director',
    '• Conceived, architected and built Civic Graph multi-tenant data warehouse processing government data from Census, Bureau of Labor Statistics, National Council of Educational Statistics',
    '• Built multi-dimensional data warehouse measuring socio-economic changes in America at every level across attitudinal, behavioral, demographic, economic and geographical dimensions',
    '• Managed engineering teams of 7-11 professionals while setting technical direction for data architecture'
) 'Heading3'

Insert-Block '• Managed engineering teams of 7-11 professionals while setting technical direction for data architecture' @(
    'Analytics Supervisor - GSD&M (Austin, TX) | November 2019 - June 2020',
    'Big Data Engineering Transformation',
    '• Transformed small data team into big data engineering team, scaling from laptop datasets to Hadoop Clusters and Hive on AWS',
    '• Managed accounts including United States Air Force, Southwest Airlines/Chase and Indeed',
    '• Rewrote mission and offerings of department and drafted integration plan with strategy team'
) 'Heading3'

Insert-Block '• Rewrote mission and offerings of department and drafted integration plan with strategy team' @(
    'Software Engineer - Mautinoa Technologies (Austin, TX) | August 2016 - February 2018',
    'SimCrisis Product Owner/Engineer',
    '• Conceived, architected and engineered econometric simulation software for humanitarian crises intervention measurement',
    '• Built SimCrisis GeoDjango web application using multi-agent modeling to create econometric simulations of crisis economies',
    '• Designed modular application accepting rules extensions for ethnic strife, different crises/disasters, supply failures'
) 'Heading3'

Insert-Block '• Designed modular application accepting rules extensions for ethnic strife, different crises/disasters, supply failures' @(
    'Senior Analyst - Myers Research (Austin, TX) | August 2012 - February 2014',
    'RACSO Product Owner/Engineer',
    '• Designed comprehensive survey instruments for specialized voting segments and niche markets',
    '• Co-developed RACSO web application managing all aspects of survey operations from instrument design to data analysis',
    '• Wrote RFP and analyzed bids from 1,200 vendors for research platform development'
) 'Heading3'

Insert-Block '• Wrote RFP and analyzed bids from 1,200 vendors for research platform development' @(
    'Research Director - PCCC (Washington, DC) | 2010 - 2012',
    'Political Research & Data Analysis (FLEEM System)',
    '• Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of simultaneous phone calls using emulated predictive dialer for regulated political surveys',
    '• Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren',
    '• Built comprehensive tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver'
) 'Heading3'

Insert-Block '• Built comprehensive tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver' @(
    'Software Engineer - Salsa Labs (Washington, DC) | January 2011 - August 2011',
    'Geospatial CRM Development',
    '• Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands simultaneously',
    '• Developed custom tile server for Web Map Service (WMS) integration using GeoTools and OpenLayers',
    '• Built advanced geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill'
) 'Heading3'

Insert-Block '• Built advanced geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill' @(
    'Programmer - Lake Research Partners (Washington, DC) | April 2008 - December 2008',
    'Political Research & Analytics',
    '• Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party',
    '• Harmonized data from 20+ polling firms with incompatible methodologies and encoding systems',
    '• Created comprehensive meta-analysis framework handling millions of survey responses that became the $400M Polling Consortium Database at The Analyst Institute, now valued at $1B+'
) 'Heading3'

# ---------------------------------------------------------------------------
# 6. Key projects - rewrite the first project in place, add two more.
# ---------------------------------------------------------------------------

Replace-ParaText 'Polling Consortium Dataset Meta-Analysis (2013 - 2016)' 'National Redistricting Platform (2020 - 2021)'
Replace-ParaText 'Comprehensive meta-analysis of polling data from tens of polling and mail firms with different methodologies and encoding systems, creating unified analytical framework' 'Cloud-based GeoDjango platform for redistricting analysis with real-time collaborative editing and Census integration, used by thousands of analysts nationwide'
Replace-ParaText 'Technologies: Python, R, Statistical Analysis, Meta-Analysis, Data Standardization' 'Technologies: GeoDjango, PostGIS, AWS, Docker, React, Python'
Replace-ParaText 'Impact: Created $400M dataset that became foundation for modern electoral analytics, estimated current value exceeds $1B' 'Impact: Reduced mapping costs by 73.5%, saving organizations $4.7M in operational expenses'

Insert-Block 'Impact: Reduced mapping costs by 73.5%, saving organizations $4.7M in operational expenses' @(
    'FLEEM Political Polling System (2010 - 2012)',
    'Completely self-built IVR system using Twilio API that contacted tens of thousands of voters daily, replicated call center functionality to performance parity',
    'Technologies: Twilio API, Python, Django, PostgreSQL, JavaScript',
    'Impact: Saved $840K in operational costs plus millions in avoided software licensing'
) 'Heading3'

Insert-Block 'Impact: Saved $840K in operational costs plus millions in avoided software licensing' @(
    'Geospatial Demographic Classification System (2013 - 2016)',
    'Machine learning platform that discovered systematic coding errors and improved demographic classification accuracy from 23% to 64%',
    'Technologies: Python, Scikit-learn, PostGIS, GeoPandas, TensorFlow',
    'Impact: Corrected demographic data affecting all Black and Asian-American voters nationwide'
) 'Heading3'

# ---------------------------------------------------------------------------
# 7. Key achievements and impact - collapse three curated sub-sections into
#    one "Impact" block with four bullets.
# ---------------------------------------------------------------------------

Delete-ParaRange 'Data Discovery' '• Interfaced with Government and Activism APIs for seamless data integration'

Insert-Block 'KEY ACHIEVEMENTS AND IMPACT' @(
    'Impact',
    '• Discovered systematic race coding errors affecting all Black and Asian-American voters',
    '• Algorithm reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M',
    '• Built redistricting platform used by thousands of analysts nationwide',
    '• Achieved 87% prediction accuracy for voter turnout vs. industry standard of 71%'
) 'Heading3'

# ---------------------------------------------------------------------------
# 8. Technical skills - keep the heading only, drop the skill lines and the
#    closing "visit my LinkedIn" paragraph.
# ---------------------------------------------------------------------------

Delete-ParaRange 'CODE Python; R; SQL; JavaScript; PHP' 'For a more detailed, full description of my experience, please visit my LinkedIn and Personal Site.'

Write-Output "DONE. Final paragraph count: $($d.Paragraphs.Count)"
